$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.453.25"
$ws.Range("E2").Value = "  -2.05%  "
$ws.Range("D3").Value = "1.843.61"
$ws.Range("E3").Value = "  -1.82%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'260.49"
$ws.Range("E5").Value = "  -7.39%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "'0.5246"
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("E8").Value = "  -8.37%  "
$ws.Range("D9").Value = "'0.06748"
$ws.Range("E9").Value = "  -4.06%  "
$ws.Range("D10").Value = "'18.92"
$ws.Range("E10").Value = "  -6.77%  "
$ws.Range("D11").Value = "'0.7716"
$ws.Range("E11").Value = "  -5.22%  "
$ws.Range("E12").Value = "  -1.25%  "
$ws.Range("D13").Value = "1.864.05"
$ws.Range("E13").Value = "  -0.72%  "
$ws.Range("D14").Value = "'89.20"
$ws.Range("E14").Value = "  -1.31%  "
$ws.Range("D15").Value = "'5.027"
$ws.Range("E15").Value = "  -3.63%  "
$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("E17").Value = "  -3.08%  "
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").Value = "'0.000007886"
$ws.Range("E19").Value = "  -3.59%  "
$ws.Range("D20").Value = "26.501.53"
$ws.Range("E20").Value = "  -1.99%  "
$ws.Range("D21").Value = "2.088.16"
$ws.Range("E21").Value = "  -1.39%  "
$ws.Range("D22").Value = "'4.536"
$ws.Range("E22").Value = "  -4.83%  "
$ws.Range("D23").Value = "'9.474"
$ws.Range("E23").Value = "  -6.99%  "
$ws.Range("D24").Value = "'5.918"
$ws.Range("E24").Value = "  -4.93%  "
$ws.Range("D25").Value = "'2.339"
$ws.Range("E25").Value = "  -2.14%  "
$ws.Range("D26").Value = "'144.26"
$ws.Range("E26").Value = "  -1.63%  "
$ws.Range("D27").Value = "'1.647"
$ws.Range("E27").Value = "  -1.67%  "
$ws.Range("D28").Value = "'16.83"
$ws.Range("E28").Value = "  -4.18%  "
$ws.Range("D29").Value = "'111.29"
$ws.Range("E29").Value = "  -1.75%  "
$ws.Range("D30").Value = "'4.190"
$ws.Range("E30").Value = "  -4.27%  "
$ws.Range("E31").Value = "  -1.31%  "
$ws.Range("E32").Value = "  -5.74%  "
$ws.Range("D33").Value = "'0.04844"
$ws.Range("E33").Value = "  -1.05%  "
$ws.Range("D34").Value = "'1.134"
$ws.Range("E34").Value = "  -3.20%  "
$ws.Range("D35").Value = "'2.846"
$ws.Range("E35").Value = "  -1.01%  "
$ws.Range("D36").Value = "'0.6863"
$ws.Range("E36").Value = "  -7.45%  "
$ws.Range("D37").Value = "'3.114"
$ws.Range("E37").Value = "  -5.57%  "
$ws.Range("D38").Value = "'0.01788"
$ws.Range("E38").Value = "  -4.99%  "
$ws.Range("D39").Value = "'2.220"
$ws.Range("E39").Value = "  -7.79%  "
$ws.Range("D40").Value = "'0.4923"
$ws.Range("E40").Value = "  -7.03%  "
$ws.Range("D41").Value = "'112.92"
$ws.Range("E41").Value = "  -3.48%  "
$ws.Range("D42").Value = "'0.8985"
$ws.Range("E42").Value = "  -8.40%  "
$ws.Range("D43").Value = "'6.173"
$ws.Range("E43").Value = "  -2.23%  "
$ws.Range("D44").Value = "'0.9998"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "'7.756"
$ws.Range("E45").Value = "  -5.15%  "
$ws.Range("D46").Value = "'0.4205"
$ws.Range("E46").Value = "  -8.72%  "
$ws.Range("D47").Value = "'0.1260"
$ws.Range("E47").Value = "  -7.59%  "
$ws.Range("D48").Value = "'9.063"
$ws.Range("E48").Value = "  -4.38%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "'35.50"
$ws.Range("E49").Value = "  -3.32%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.05877"
$ws.Range("E50").Value = "  -1.18%  "
$ws.Range("D51").Value = "'59.30"
$ws.Range("E51").Value = "  -4.09%  "

# Reset style on text-protected numeric-looking cells to clear quote-prefix formatting
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
